$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for all data rows (2-118) from 45203 to 45205
$ws.Range("C2:C118").Value = 45205
